$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 5 with trade data, matching the style of existing rows
$ws.Range("A5").Value = 9959.23
$ws.Range("B5").Value = 10035.5
$ws.Range("C5").Value = 109.08
$ws.Range("D5").Value = 108.25
$ws.Range("E5").Value = $false
$ws.Range("F5").Value = -0.76
$ws.Range("G5").Value = 42612.674525462964
$ws.Range("H5").Value = $false

# Match the date/time cell style used by the other cells in column G
# (copy format only, so we reuse the existing style instead of creating a new one)
$ws.Range("G4").Copy()
$ws.Range("G5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
